$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "2+32="
$t.Cell(1, 2).Range.Text = "89+6="
$t.Cell(1, 3).Range.Text = "6+1="
$t.Cell(1, 4).Range.Text = "94-39="
$t.Cell(1, 5).Range.Text = "9-0="
$t.Cell(2, 1).Range.Text = "19+65="
$t.Cell(2, 2).Range.Text = "33+37="
$t.Cell(2, 3).Range.Text = "34+48="
$t.Cell(2, 4).Range.Text = "48-19="
$t.Cell(2, 5).Range.Text = "84-66="
$t.Cell(3, 1).Range.Text = "31-13="
$t.Cell(3, 2).Range.Text = "85+14="
$t.Cell(3, 3).Range.Text = "65-60="
$t.Cell(3, 4).Range.Text = "0+1="
$t.Cell(3, 5).Range.Text = "4+27="
$t.Cell(4, 1).Range.Text = "49-37="
$t.Cell(4, 2).Range.Text = "9+7="
$t.Cell(4, 3).Range.Text = "82+12="
$t.Cell(4, 4).Range.Text = "17+56="
$t.Cell(4, 5).Range.Text = "28+28="
$t.Cell(5, 1).Range.Text = "24-8="
$t.Cell(5, 2).Range.Text = "90-82="
$t.Cell(5, 3).Range.Text = "62-27="
$t.Cell(5, 4).Range.Text = "76-25="
$t.Cell(5, 5).Range.Text = "11+66="
$t.Cell(6, 1).Range.Text = "36+8="
$t.Cell(6, 2).Range.Text = "55+0="
$t.Cell(6, 3).Range.Text = "56+7="
$t.Cell(6, 4).Range.Text = "90-57="
$t.Cell(6, 5).Range.Text = "46-27="
$t.Cell(7, 1).Range.Text = "5+82="
$t.Cell(7, 2).Range.Text = "50+20="
$t.Cell(7, 3).Range.Text = "71-49="
$t.Cell(7, 4).Range.Text = "97+2="
$t.Cell(7, 5).Range.Text = "11+43="
$t.Cell(8, 1).Range.Text = "89-16="
$t.Cell(8, 2).Range.Text = "19+17="
$t.Cell(8, 3).Range.Text = "50+36="
$t.Cell(8, 4).Range.Text = "24+40="
$t.Cell(8, 5).Range.Text = "76-42="
$t.Cell(9, 1).Range.Text = "89-5="
$t.Cell(9, 2).Range.Text = "30+11="
$t.Cell(9, 3).Range.Text = "36+43="
$t.Cell(9, 4).Range.Text = "97-15="
$t.Cell(9, 5).Range.Text = "96-58="
$t.Cell(10, 1).Range.Text = "27+30="
$t.Cell(10, 2).Range.Text = "80-20="
$t.Cell(10, 3).Range.Text = "85-1="
$t.Cell(10, 4).Range.Text = "39+37="
$t.Cell(10, 5).Range.Text = "83-3="
$t.Cell(11, 1).Range.Text = "2+63="
$t.Cell(11, 2).Range.Text = "4+84="
$t.Cell(11, 3).Range.Text = "65-2="
$t.Cell(11, 4).Range.Text = "83-12="
$t.Cell(11, 5).Range.Text = "70-3="
$t.Cell(12, 1).Range.Text = "32+0="
$t.Cell(12, 2).Range.Text = "39+0="
$t.Cell(12, 3).Range.Text = "49-4="
$t.Cell(12, 4).Range.Text = "94-91="
$t.Cell(12, 5).Range.Text = "84-52="
$t.Cell(13, 1).Range.Text = "76-4="
$t.Cell(13, 2).Range.Text = "68-28="
$t.Cell(13, 3).Range.Text = "74+4="
$t.Cell(13, 4).Range.Text = "91-68="
$t.Cell(13, 5).Range.Text = "22+45="
$t.Cell(14, 1).Range.Text = "82-10="
$t.Cell(14, 2).Range.Text = "15+79="
$t.Cell(14, 3).Range.Text = "90-82="
$t.Cell(14, 4).Range.Text = "82-61="
$t.Cell(14, 5).Range.Text = "77+12="
$t.Cell(15, 1).Range.Text = "52-1="
$t.Cell(15, 2).Range.Text = "5+28="
$t.Cell(15, 3).Range.Text = "34-10="
$t.Cell(15, 4).Range.Text = "86-4="
$t.Cell(15, 5).Range.Text = "32+52="
$t.Cell(16, 1).Range.Text = "22+61="
$t.Cell(16, 2).Range.Text = "48-26="
$t.Cell(16, 3).Range.Text = "3+68="
$t.Cell(16, 4).Range.Text = "66-4="
$t.Cell(16, 5).Range.Text = "40+50="
$t.Cell(17, 1).Range.Text = "7+8="
$t.Cell(17, 2).Range.Text = "63-12="
$t.Cell(17, 3).Range.Text = "78-34="
$t.Cell(17, 4).Range.Text = "71-31="
$t.Cell(17, 5).Range.Text = "57-45="
$t.Cell(18, 1).Range.Text = "6+15="
$t.Cell(18, 2).Range.Text = "78-42="
$t.Cell(18, 3).Range.Text = "14+55="
$t.Cell(18, 4).Range.Text = "99-63="
$t.Cell(18, 5).Range.Text = "74-16="
$t.Cell(19, 1).Range.Text = "69-26="
$t.Cell(19, 2).Range.Text = "67+4="
$t.Cell(19, 3).Range.Text = "3+40="
$t.Cell(19, 4).Range.Text = "15-12="
$t.Cell(19, 5).Range.Text = "55+33="
$t.Cell(20, 1).Range.Text = "44+44="
$t.Cell(20, 2).Range.Text = "93-86="
$t.Cell(20, 3).Range.Text = "99-49="
$t.Cell(20, 4).Range.Text = "38+57="
$t.Cell(20, 5).Range.Text = "28+47="
